# Add two new columns "I0" (I) and "IF" (J) to the right of the existing
# data table (which currently ends at column H), mirroring the header
# style used by the other header cells, then fill in the per-row values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new labels "I0" / "IF", matching the style of
#     the existing bold/bordered header cells (copy format from H1). ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-35: columns I (I0) and J (IF) ---
$data = New-Object 'object[,]' 34,2

$data[0,0]  = 6; $data[0,1]  = 6    # row 2
$data[1,0]  = 5; $data[1,1]  = 6    # row 3
$data[2,0]  = 9; $data[2,1]  = 9    # row 4
$data[3,0]  = 8; $data[3,1]  = 8    # row 5
$data[4,0]  = 8; $data[4,1]  = 8    # row 6
$data[5,0]  = 9; $data[5,1]  = 9    # row 7
$data[6,0]  = 7; $data[6,1]  = 8    # row 8
$data[7,0]  = 6; $data[7,1]  = 6    # row 9
$data[8,0]  = 7; $data[8,1]  = 7    # row 10
$data[9,0]  = 7; $data[9,1]  = 7    # row 11
$data[10,0] = 9; $data[10,1] = 9    # row 12
$data[11,0] = 6; $data[11,1] = 7    # row 13
$data[12,0] = 8; $data[12,1] = 8    # row 14
$data[13,0] = 7; $data[13,1] = 7    # row 15
$data[14,0] = 4; $data[14,1] = 5    # row 16
$data[15,0] = 7; $data[15,1] = 7    # row 17
$data[16,0] = 8; $data[16,1] = 8    # row 18
$data[17,0] = 6; $data[17,1] = 7    # row 19
$data[18,0] = 7; $data[18,1] = 7    # row 20
$data[19,0] = 9; $data[19,1] = 9    # row 21
$data[20,0] = 9; $data[20,1] = 9    # row 22
$data[21,0] = 7; $data[21,1] = 7    # row 23
$data[22,0] = 6; $data[22,1] = 6    # row 24
$data[23,0] = 6; $data[23,1] = 6    # row 25
$data[24,0] = 6; $data[24,1] = 7    # row 26
$data[25,0] = 6; $data[25,1] = 7    # row 27
$data[26,0] = 5; $data[26,1] = 5    # row 28
$data[27,0] = 7; $data[27,1] = 7    # row 29
$data[28,0] = 5; $data[28,1] = 6    # row 30
$data[29,0] = 8; $data[29,1] = 9    # row 31
$data[30,0] = 5; $data[30,1] = 5    # row 32
$data[31,0] = 9; $data[31,1] = 9    # row 33
$data[32,0] = 7; $data[32,1] = 7    # row 34
$data[33,0] = 8; $data[33,1] = 8    # row 35

$ws.Range("I2:J35").Value = $data
